$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.614.84'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.482.90'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9703'
$ws.Range("E5").Value = '  +2.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '280.50'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3669'
$ws.Range("E7").Value = '  -1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.19'
$ws.Range("E9").Value = '  -2.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.068'
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06697'
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.546'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.23'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.232'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9706'
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.483.32'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05971'
$ws.Range("E19").Value = '  +4.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.00'
$ws.Range("E20").Value = '  -3.29%  '
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.275'
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.664.21'
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.56'
$ws.Range("E26").Value = '  +3.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.134'
$ws.Range("E27").Value = '  -7.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.36'
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.644.56'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.82'
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.946'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8280'
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.053'
$ws.Range("E33").Value = '  -5.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08035'
$ws.Range("E34").Value = '  +2.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.540'
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.213'
$ws.Range("E36").Value = '  +7.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05816'
$ws.Range("E37").Value = '  -3.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.785'
$ws.Range("E38").Value = '  -3.04%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.712'
$ws.Range("E39").Value = '  +1.91%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02053'
$ws.Range("E40").Value = '  -0.98%  '
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9702'
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1889'
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5332'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.548'
$ws.Range("E45").Value = '  -1.18%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.28'
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.25'
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5227'
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06513'
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9947'
$ws.Range("E51").Value = '  +0.27%  '
